$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Mmrn2"
$ws.Range("C2").Value = "Cd93"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 27.10983833333333
$ws.Range("H2").Value = 81.329515
$ws.Range("I2").Value = 0.9284397459331688
$ws.Range("J2").Value = 0.9284397459331687
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 124.722578
$ws.Range("N2").Value = 374.167734
$ws.Range("O2").Value = 0.9767311432246923
$ws.Range("P2").Value = 0.9767311432246923
$ws.Range("Q2").Value = 3381.208926096556
$ws.Range("R2").Value = 30430.88033486901
$ws.Range("S2").Value = 0.9068360144605468
$ws.Range("T2").Value = 0.9068360144605467

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Mmrn2"
$ws.Range("C3").Value = "Cd93"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 27.10983833333333
$ws.Range("H3").Value = 81.329515
$ws.Range("I3").Value = 0.9284397459331688
$ws.Range("J3").Value = 0.9284397459331687
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.134712
$ws.Range("N3").Value = 0.404136
$ws.Range("O3").Value = 0.001054960600366076
$ws.Range("P3").Value = 0.001054960600366076
$ws.Range("Q3").Value = 3.65202054156
$ws.Range("R3").Value = 32.86818487404
$ws.Range("S3").Value = 0.0009794673517733829
$ws.Range("T3").Value = 0.0009794673517733827

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Mmrn2"
$ws.Range("C4").Value = "Cd93"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 27.10983833333333
$ws.Range("H4").Value = 81.329515
$ws.Range("I4").Value = 0.9284397459331688
$ws.Range("J4").Value = 0.9284397459331687
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.836578333333333
$ws.Range("N4").Value = 8.509734999999999
$ws.Range("O4").Value = 0.02221389617494163
$ws.Range("P4").Value = 0.02221389617494163
$ws.Range("Q4").Value = 76.89918003650277
$ws.Range("R4").Value = 692.092620328525
$ws.Range("S4").Value = 0.02062426412084859
$ws.Range("T4").Value = 0.02062426412084859

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Mmrn2"
$ws.Range("C5").Value = "Cd93"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.190813333333333
$ws.Range("H5").Value = 3.57244
$ws.Range("I5").Value = 0.04078218449921273
$ws.Range("J5").Value = 0.04078218449921273
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 124.722578
$ws.Range("N5").Value = 374.167734
$ws.Range("O5").Value = 0.9767311432246923
$ws.Range("P5").Value = 0.9767311432246923
$ws.Range("Q5").Value = 148.5213088501067
$ws.Range("R5").Value = 1336.69177965096
$ws.Range("S5").Value = 0.03983322968911638
$ws.Range("T5").Value = 0.03983322968911638

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Mmrn2"
$ws.Range("C6").Value = "Cd93"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.190813333333333
$ws.Range("H6").Value = 3.57244
$ws.Range("I6").Value = 0.04078218449921273
$ws.Range("J6").Value = 0.04078218449921273
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.134712
$ws.Range("N6").Value = 0.404136
$ws.Range("O6").Value = 0.001054960600366076
$ws.Range("P6").Value = 0.001054960600366076
$ws.Range("Q6").Value = 0.16041684576
$ws.Range("R6").Value = 1.44375161184
$ws.Range("S6").Value = 0.00004302359784352954
$ws.Range("T6").Value = 0.00004302359784352954

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Mmrn2"
$ws.Range("C7").Value = "Cd93"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.190813333333333
$ws.Range("H7").Value = 3.57244
$ws.Range("I7").Value = 0.04078218449921273
$ws.Range("J7").Value = 0.04078218449921273
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.836578333333333
$ws.Range("N7").Value = 8.509734999999999
$ws.Range("O7").Value = 0.02221389617494163
$ws.Range("P7").Value = 0.02221389617494163
$ws.Range("Q7").Value = 3.377835300377777
$ws.Range("R7").Value = 30.40051770339999
$ws.Range("S7").Value = 0.0009059312122528255
$ws.Range("T7").Value = 0.0009059312122528255

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Mmrn2"
$ws.Range("C8").Value = "Cd93"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.8986996666666666
$ws.Range("H8").Value = 2.696099
$ws.Range("I8").Value = 0.03077806956761847
$ws.Range("J8").Value = 0.03077806956761847
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 124.722578
$ws.Range("N8").Value = 374.167734
$ws.Range("O8").Value = 0.9767311432246923
$ws.Range("P8").Value = 0.9767311432246923
$ws.Range("Q8").Value = 112.0881392744073
$ws.Range("R8").Value = 1008.793253469666
$ws.Range("S8").Value = 0.0300618990750291
$ws.Range("T8").Value = 0.0300618990750291

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Mmrn2"
$ws.Range("C9").Value = "Cd93"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.8986996666666666
$ws.Range("H9").Value = 2.696099
$ws.Range("I9").Value = 0.03077806956761847
$ws.Range("J9").Value = 0.03077806956761847
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.134712
$ws.Range("N9").Value = 0.404136
$ws.Range("O9").Value = 0.001054960600366076
$ws.Range("P9").Value = 0.001054960600366076
$ws.Range("Q9").Value = 0.121065629496
$ws.Range("R9").Value = 1.089590665464
$ws.Range("S9").Value = 0.00003246965074916364
$ws.Range("T9").Value = 0.00003246965074916364

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Mmrn2"
$ws.Range("C10").Value = "Cd93"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.8986996666666666
$ws.Range("H10").Value = 2.696099
$ws.Range("I10").Value = 0.03077806956761847
$ws.Range("J10").Value = 0.03077806956761847
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.836578333333333
$ws.Range("N10").Value = 8.509734999999999
$ws.Range("O10").Value = 0.02221389617494163
$ws.Range("P10").Value = 0.02221389617494163
$ws.Range("Q10").Value = 2.549232002640555
$ws.Range("R10").Value = 22.943088023765
$ws.Range("S10").Value = 0.0006837008418402073
$ws.Range("T10").Value = 0.0006837008418402073
